$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly data between row 2 and row 3 (Fecha, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Precio $/Kg) so the rows are in
# chronological order.

$cols = @("D", "M", "N", "O", "P", "S")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")
    $v2 = $cell2.Value2
    $v3 = $cell3.Value2
    $cell2.Value = $v3
    $cell3.Value = $v2
}
